$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.711932795781365
$ws.Range("C2").Value = 0.3813972416621709
$ws.Range("E2").Value = 0.1673490678546612
$ws.Range("F2").Value = 2.980879819969289
$ws.Range("G2").Value = 0.002496035156098531
$ws.Range("J2").Value = 0.1091512377589368
$ws.Range("M2").Value = 0.5741148005297987
$ws.Range("N2").Value = 1.864864243144368
$ws.Range("B3").Value = 1.577858537099019
$ws.Range("C3").Value = 0.3427137432935581
$ws.Range("E3").Value = 0.1674625896059725
$ws.Range("F3").Value = 2.934652647296019
$ws.Range("G3").Value = 0.002502206836288865
$ws.Range("J3").Value = 0.1086827507448263
$ws.Range("M3").Value = 0.5497237451810264
$ws.Range("N3").Value = 1.874376330090996
$ws.Range("B4").Value = 1.496583182097481
$ws.Range("C4").Value = 0.3191586110819173
$ws.Range("E4").Value = 0.1675568381340966
$ws.Range("F4").Value = 2.907977110180909
$ws.Range("G4").Value = 0.002506191323185795
$ws.Range("J4").Value = 0.1084533633033047
$ws.Range("M4").Value = 0.5351109618494689
$ws.Range("N4").Value = 1.880851983922923
$ws.Range("B5").Value = 1.463722385518224
$ws.Range("C5").Value = 0.3096079616137217
$ws.Range("E5").Value = 0.167601434749419
$ws.Range("F5").Value = 2.897534129360821
$ws.Range("G5").Value = 0.002507864262264302
$ws.Range("J5").Value = 0.1083744234450457
$ws.Range("M5").Value = 0.5292471510903312
$ws.Range("N5").Value = 1.883649715426756
$ws.Range("B6").Value = 1.458281447153013
$ws.Range("C6").Value = 0.3080249634323025
$ws.Range("E6").Value = 0.1676092143494046
$ws.Range("F6").Value = 2.895825843274864
$ws.Range("G6").Value = 0.002508145030945486
$ws.Range("J6").Value = 0.108362190257342
$ws.Range("M6").Value = 0.5282789564584931
$ws.Range("N6").Value = 1.884123845298362
$ws.Range("B7").Value = 1.496138963692204
$ws.Range("C7").Value = 0.319029613863961
$ws.Range("E7").Value = 0.1675574144943441
$ws.Range("F7").Value = 2.907834544106564
$ws.Range("G7").Value = 0.00250621368544568
$ws.Range("J7").Value = 0.1084522399806112
$ws.Range("M7").Value = 0.5350315123987457
$ws.Range("N7").Value = 1.880889073075664
$ws.Range("B8").Value = 1.665484475037488
$ws.Range("C8").Value = 0.3680176083916535
$ws.Range("E8").Value = 0.1673831256976728
$ws.Range("F8").Value = 2.964584929870838
$ws.Range("G8").Value = 0.002498122778328593
$ws.Range("J8").Value = 0.1089775411179517
$ws.Range("M8").Value = 0.5656291020782263
$ws.Range("N8").Value = 1.868011647734903
$ws.Range("B9").Value = 2.006060988782679
$ws.Range("C9").Value = 0.4657032025084504
$ws.Range("E9").Value = 0.1672354381376975
$ws.Range("F9").Value = 3.089526401137874
$ws.Range("G9").Value = 0.002483795870423927
$ws.Range("J9").Value = 0.1104753080719831
$ws.Range("M9").Value = 0.6285354141144239
$ws.Range("N9").Value = 1.847838749994565
$ws.Range("B10").Value = 2.26174337469763
$ws.Range("C10").Value = 0.5385509760610034
$ws.Range("E10").Value = 0.1672445029719478
$ws.Range("F10").Value = 3.189799469931103
$ws.Range("G10").Value = 0.002474196658885709
$ws.Range("J10").Value = 0.111868505210154
$ws.Range("M10").Value = 0.6765585753851013
$ws.Range("N10").Value = 1.83616982434657
$ws.Range("B11").Value = 2.37930988996817
$ws.Range("C11").Value = 0.5719456956335875
$ws.Range("E11").Value = 0.1672740249501192
$ws.Range("F11").Value = 3.237292187923458
$ws.Range("G11").Value = 0.002470028470857719
$ws.Range("J11").Value = 0.1125675884751161
$ws.Range("M11").Value = 0.6988058614233097
$ws.Range("N11").Value = 1.831557931449609
$ws.Range("B12").Value = 2.424014545545617
$ws.Range("C12").Value = 0.5846297309553847
$ws.Range("E12").Value = 0.1672888460756914
$ws.Range("F12").Value = 3.255549121038655
$ws.Range("G12").Value = 0.002468478444606871
$ws.Range("J12").Value = 0.1128418424340509
$ws.Range("M12").Value = 0.7072886023451588
$ws.Range("N12").Value = 1.82991268244335
$ws.Range("B13").Value = 2.414378311820542
$ws.Range("C13").Value = 0.5818962766784921
$ws.Range("E13").Value = 0.1672854922656395
$ws.Range("F13").Value = 3.251605006573186
$ws.Range("G13").Value = 0.002468811011252736
$ws.Range("J13").Value = 0.1127823510831405
$ws.Range("M13").Value = 0.7054590972697241
$ws.Range("N13").Value = 1.83026249977128
$ws.Range("B14").Value = 2.382984041183875
$ws.Range("C14").Value = 0.5729884456672494
$ws.Range("E14").Value = 0.1672751713612932
$ws.Range("F14").Value = 3.23878872204719
$ws.Range("G14").Value = 0.002469900381406911
$ws.Range("J14").Value = 0.1125899598124462
$ws.Range("M14").Value = 0.6995025729410571
$ws.Range("N14").Value = 1.831420542130601
$ws.Range("B15").Value = 2.363778333911796
$ws.Range("C15").Value = 0.567537155541288
$ws.Range("E15").Value = 0.1672693234899167
$ws.Range("F15").Value = 3.230973934111688
$ws.Range("G15").Value = 0.002470571342489852
$ws.Range("J15").Value = 0.1124733592451435
$ws.Range("M15").Value = 0.6958616202701791
$ws.Range("N15").Value = 1.832143083529218
$ws.Range("B16").Value = 2.254085748346313
$ws.Range("C16").Value = 0.5363738195024439
$ws.Range("E16").Value = 0.167243083798847
$ws.Range("F16").Value = 3.186733704827617
$ws.Range("G16").Value = 0.00247447303996035
$ws.Range("J16").Value = 0.1118241441574952
$ws.Range("M16").Value = 0.6751127780164126
$ws.Range("N16").Value = 1.83648532342616
$ws.Range("B17").Value = 2.187117416715012
$ws.Range("C17").Value = 0.5173226396659061
$ws.Range("E17").Value = 0.1672334838918639
$ws.Range("F17").Value = 3.16007646260627
$ws.Range("G17").Value = 0.002476917332305196
$ws.Range("J17").Value = 0.1114426945216636
$ws.Range("M17").Value = 0.6624871490405297
$ws.Range("N17").Value = 1.839328236183263
$ws.Range("B18").Value = 2.148716900178101
$ws.Range("C18").Value = 0.5063889074551753
$ws.Range("E18").Value = 0.1672303534392245
$ws.Range("F18").Value = 3.144920540439216
$ws.Range("G18").Value = 0.002478341922849579
$ws.Range("J18").Value = 0.1112294339624071
$ws.Range("M18").Value = 0.6552629653585598
$ws.Range("N18").Value = 1.841028886131454
$ws.Range("B19").Value = 2.135735265699338
$ws.Range("C19").Value = 0.5026910143902228
$ws.Range("E19").Value = 0.1672297045455657
$ws.Range("F19").Value = 3.139819267067224
$ws.Range("G19").Value = 0.002478827481291205
$ws.Range("J19").Value = 0.1111582779449733
$ws.Range("M19").Value = 0.6528234494494569
$ws.Range("N19").Value = 1.8416159120858
$ws.Range("B20").Value = 2.194234077547492
$ws.Range("C20").Value = 0.5193481761852468
$ws.Range("E20").Value = 0.1672342584220416
$ws.Range("F20").Value = 3.162895874842803
$ws.Range("G20").Value = 0.00247665519945803
$ws.Range("J20").Value = 0.1114826641364957
$ws.Range("M20").Value = 0.6638272590731304
$ws.Range("N20").Value = 1.839018818679619
$ws.Range("B21").Value = 2.392200250592339
$ws.Range("C21").Value = 0.575603844239879
$ws.Range("E21").Value = 0.1672781041011842
$ws.Range("F21").Value = 3.242545761178235
$ws.Range("G21").Value = 0.002469579638340293
$ws.Range("J21").Value = 0.1126462101664885
$ws.Range("M21").Value = 0.7012505658196346
$ws.Range("N21").Value = 1.831077642810769
$ws.Range("B22").Value = 2.522661593706914
$ws.Range("C22").Value = 0.6125933759990403
$ws.Range("E22").Value = 0.1673279829847552
$ws.Range("F22").Value = 3.296190894886848
$ws.Range("G22").Value = 0.002465120667346618
$ws.Range("J22").Value = 0.1134622454581162
$ws.Range("M22").Value = 0.7260481858405967
$ws.Range("N22").Value = 1.82647784163359
$ws.Range("B23").Value = 2.452931761986576
$ws.Range("C23").Value = 0.5928304806996039
$ws.Range("E23").Value = 0.1672994228935245
$ws.Range("F23").Value = 3.267413242473538
$ws.Range("G23").Value = 0.002467485434936286
$ws.Range("J23").Value = 0.1130215804175947
$ws.Range("M23").Value = 0.7127820221915044
$ws.Range("N23").Value = 1.828878490162666
$ws.Range("B24").Value = 2.191016322565872
$ws.Range("C24").Value = 0.518432371908375
$ws.Range("E24").Value = 0.1672339008159671
$ws.Range("F24").Value = 3.161620690232809
$ws.Range("G24").Value = 0.002476773649556042
$ws.Range("J24").Value = 0.1114645750708334
$ws.Range("M24").Value = 0.6632212880316501
$ws.Range("N24").Value = 1.839158500146297
$ws.Range("B25").Value = 1.91298790455528
$ws.Range("C25").Value = 0.4390939395852342
$ws.Range("E25").Value = 0.1672546947626694
$ws.Range("F25").Value = 3.054248825123835
$ws.Range("G25").Value = 0.002487508088267019
$ws.Range("J25").Value = 0.1100192691878945
$ws.Range("M25").Value = 0.6112030510931419
$ws.Range("N25").Value = 1.852746560176257
